$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3424.2114
$ws.Range("I132").Value = 3091.5945
$ws.Range("J132").Value = 4244.6665
$ws.Range("K132").Value = 9274.783500000001
$ws.Range("L132").Value = 12733.9995
$ws.Range("M132").Value = -6744.783500000001
$ws.Range("N132").Value = -17793.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 863.53845
$ws.Range("I2").Value = 838.6
$ws.Range("J2").Value = 946.6667
$ws.Range("K2").Value = 838.6
$ws.Range("L2").Value = 946.6667
$ws.Range("M2").Value = -725.6
$ws.Range("N2").Value = -1172.6667

$ws.Range("H4").Value = 254.90909
$ws.Range("I4").Value = 166.66667
$ws.Range("J4").Value = 360.8
$ws.Range("K4").Value = 166.66667
$ws.Range("L4").Value = 360.8
$ws.Range("M4").Value = -50.66667000000001
$ws.Range("N4").Value = -592.8

$ws.Range("H6").Value = 9000.4
$ws.Range("I6").Value = 10000.5
$ws.Range("K6").Value = 10000.5
$ws.Range("M6").Value = -9827.5

$ws.Range("H16").Value = 3002
$ws.Range("I16").Value = 3002
$ws.Range("K16").Value = 3002
$ws.Range("M16").Value = -2715

$ws.Range("H32").Value = 11420.107
$ws.Range("I32").Value = 3059.7273
$ws.Range("K32").Value = 3059.7273
$ws.Range("M32").Value = -2772.7273

$ws.Range("H61").Value = 810.75
$ws.Range("I61").Value = 652.7368
$ws.Range("J61").Value = 1144.3334
$ws.Range("K61").Value = 652.7368
$ws.Range("L61").Value = 1144.3334
$ws.Range("M61").Value = -440.7368
$ws.Range("N61").Value = -1568.3334

$ws.Range("H88").Value = 66722148
$ws.Range("I88").Value = 1839.8
$ws.Range("J88").Value = 100082296
$ws.Range("K88").Value = 1839.8
$ws.Range("L88").Value = 100082296
$ws.Range("M88").Value = -1433.8
$ws.Range("N88").Value = -100083108

$ws.Range("H91").Value = 66722148
$ws.Range("I91").Value = 1839.8
$ws.Range("J91").Value = 100082296
$ws.Range("K91").Value = 1839.8
$ws.Range("L91").Value = 100082296
$ws.Range("M91").Value = -435.8
$ws.Range("N91").Value = -100085104

$ws.Range("H106").Value = 43750
$ws.Range("J106").Value = 43750
$ws.Range("L106").Value = 43750
$ws.Range("N106").Value = -46274

$ws.Range("H116").Value = 863.53845
$ws.Range("I116").Value = 838.6
$ws.Range("J116").Value = 946.6667
$ws.Range("K116").Value = 838.6
$ws.Range("L116").Value = 946.6667
$ws.Range("M116").Value = 1455.4
$ws.Range("N116").Value = -5534.6667

$ws.Range("H122").Value = 2425.24
$ws.Range("I122").Value = 1839.1111
$ws.Range("J122").Value = 2754.9375
$ws.Range("K122").Value = 5517.3333
$ws.Range("L122").Value = 8264.8125
$ws.Range("M122").Value = -3067.3333
$ws.Range("N122").Value = -13164.8125

$ws.Range("H136").Value = 810.75
$ws.Range("I136").Value = 652.7368
$ws.Range("J136").Value = 1144.3334
$ws.Range("K136").Value = 1958.2104
$ws.Range("L136").Value = 3433.0002
$ws.Range("M136").Value = 591.7896000000001
$ws.Range("N136").Value = -8533.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 863.53845
$ws.Range("I3").Value = 838.6
$ws.Range("J3").Value = 946.6667
$ws.Range("K3").Value = 838.6
$ws.Range("L3").Value = 946.6667
$ws.Range("M3").Value = -724.6
$ws.Range("N3").Value = -1174.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 33336834
$ws.Range("I16").Value = 41670000
$ws.Range("K16").Value = 41670000
$ws.Range("M16").Value = -41669713

$ws.Range("H31").Value = 6932524.5
$ws.Range("I31").Value = 8740178
$ws.Range("J31").Value = 3185.6667
$ws.Range("K31").Value = 8740178
$ws.Range("L31").Value = 3185.6667
$ws.Range("M31").Value = -8739883
$ws.Range("N31").Value = -3775.6667

$ws.Range("H34").Value = 6932524.5
$ws.Range("I34").Value = 8740178
$ws.Range("J34").Value = 3185.6667
$ws.Range("K34").Value = 8740178
$ws.Range("L34").Value = 3185.6667
$ws.Range("M34").Value = -8739976
$ws.Range("N34").Value = -3589.6667

$ws.Range("H58").Value = 791.9836
$ws.Range("I58").Value = 478.18604
$ws.Range("J58").Value = 1541.6111
$ws.Range("K58").Value = 478.18604
$ws.Range("L58").Value = 1541.6111
$ws.Range("M58").Value = -275.18604
$ws.Range("N58").Value = -1947.6111

$ws.Range("H62").Value = 58826296
$ws.Range("I62").Value = 2925.25
$ws.Range("J62").Value = 76925790
$ws.Range("K62").Value = 2925.25
$ws.Range("L62").Value = 76925790
$ws.Range("M62").Value = -2301.25
$ws.Range("N62").Value = -76927038

$ws.Range("H65").Value = 58826296
$ws.Range("I65").Value = 2925.25
$ws.Range("J65").Value = 76925790
$ws.Range("K65").Value = 14626.25
$ws.Range("L65").Value = 384628950
$ws.Range("M65").Value = -11506.25
$ws.Range("N65").Value = -384635190

$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996

$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984

$ws.Range("H105").Value = 844.2857
$ws.Range("I105").Value = 577.7778
$ws.Range("J105").Value = 1324
$ws.Range("K105").Value = 577.7778
$ws.Range("L105").Value = 1324
$ws.Range("M105").Value = 1169.2222
$ws.Range("N105").Value = -4818

$ws.Range("H107").Value = 22728234
$ws.Range("I107").Value = 38462500
$ws.Range("J107").Value = 960.7778
$ws.Range("K107").Value = 38462500
$ws.Range("L107").Value = 960.7778
$ws.Range("M107").Value = -38460580
$ws.Range("N107").Value = -4800.7778

$ws.Range("H113").Value = 33336834
$ws.Range("I113").Value = 41670000
$ws.Range("K113").Value = 41670000
$ws.Range("M113").Value = -41667830

$ws.Range("H136").Value = 791.9836
$ws.Range("I136").Value = 478.18604
$ws.Range("J136").Value = 1541.6111
$ws.Range("K136").Value = 1434.55812
$ws.Range("L136").Value = 4624.8333
$ws.Range("M136").Value = 1115.44188
$ws.Range("N136").Value = -9724.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 412.5
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 425
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 1275
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -3771

$ws.Range("H98").Value = 408.10526
$ws.Range("J98").Value = 385.2
$ws.Range("L98").Value = 1155.6
$ws.Range("N98").Value = -4151.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 28269526
$ws.Range("I68").Value = 112778450
$ws.Range("J68").Value = 99886
$ws.Range("K68").Value = 112778450
$ws.Range("L68").Value = 99886
$ws.Range("M68").Value = -112777701
$ws.Range("N68").Value = -101384

$ws.Range("H71").Value = 28269526
$ws.Range("I71").Value = 112778450
$ws.Range("J71").Value = 99886
$ws.Range("K71").Value = 563892250
$ws.Range("L71").Value = 499430
$ws.Range("M71").Value = -563888506
$ws.Range("N71").Value = -506918

$ws.Range("H108").Value = 53500
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 53500
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 53500
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -61180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H94").Value = 28650
$ws.Range("J94").Value = 28650
$ws.Range("L94").Value = 28650
$ws.Range("N94").Value = -30452

$ws.Range("H104").Value = 27275
$ws.Range("J104").Value = 27275
$ws.Range("L104").Value = 27275
$ws.Range("N104").Value = -34263
